# Auto-generated Excel COM-interop script applying the cryptos.xlsx data refresh
# described by the commit 'Updated cryptos list on Wed Jul 19 19:58:14 UTC 2023 with GitHub Actions'.

function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '30.089.71'
Set-TextValue $ws 'E2' '  +1.08%  '
Set-TextValue $ws 'D3' '1.912.99'
Set-TextValue $ws 'E3' '  +0.96%  '
Set-TextValue $ws 'D4' '1.000'
Set-TextValue $ws 'E4' '  -0.13%  '
Set-TextValue $ws 'D5' '0.8398'
Set-TextValue $ws 'D6' '242.36'
Set-TextValue $ws 'E6' '  +0.81%  '
Set-TextValue $ws 'D7' '1.001'
Set-TextValue $ws 'E7' '  -0.02%  '
Set-TextValue $ws 'D8' '0.3251'
Set-TextValue $ws 'E8' '  +6.81%  '
Set-TextValue $ws 'D9' '26.86'
Set-TextValue $ws 'E9' '  +5.94%  '
Set-TextValue $ws 'D10' '0.07110'
Set-TextValue $ws 'E10' '  +4.43%  '
Set-TextValue $ws 'D11' '0.08050'
Set-TextValue $ws 'E11' '  +1.10%  '
Set-TextValue $ws 'D12' '0.7561'
Set-TextValue $ws 'E12' '  +2.83%  '
Set-TextValue $ws 'D13' '1.910.90'
Set-TextValue $ws 'E13' '  +0.90%  '
Set-TextValue $ws 'D14' '5.232'
Set-TextValue $ws 'E14' '  +1.83%  '
Set-TextValue $ws 'D15' '93.13'
Set-TextValue $ws 'E15' '  +2.66%  '
Set-TextValue $ws 'D16' '14.25'
Set-TextValue $ws 'E16' '  +3.13%  '
Set-TextValue $ws 'D17' '30.086.94'
Set-TextValue $ws 'E17' '  +1.05%  '
Set-TextValue $ws 'D18' '5.975'
Set-TextValue $ws 'E18' '  +1.20%  '
Set-TextValue $ws 'D19' '245.57'
Set-TextValue $ws 'E19' '  +1.63%  '
Set-TextValue $ws 'D20' '0.000007801'
Set-TextValue $ws 'E20' '  +1.63%  '
Set-TextValue $ws 'B21' 'Dai'
Set-TextValue $ws 'C21' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D21' '1.001'
Set-TextValue $ws 'E21' '  +0.01%  '
Set-TextValue $ws 'B22' 'WrappedliquidstakedEther2.0'
Set-TextValue $ws 'C22' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws 'D22' '2.155.00'
Set-TextValue $ws 'E22' '  +0.64%  '
Set-TextValue $ws 'D23' '1.000'
Set-TextValue $ws 'E23' '  -0.11%  '
Set-TextValue $ws 'D24' '7.047'
Set-TextValue $ws 'E24' '  +2.27%  '
Set-TextValue $ws 'D25' '0.1606'
Set-TextValue $ws 'E25' '  +24.71%  '
Set-TextValue $ws 'D26' '169.39'
Set-TextValue $ws 'E26' '  +1.89%  '
Set-TextValue $ws 'D27' '9.304'
Set-TextValue $ws 'E27' '  +1.38%  '
Set-TextValue $ws 'D28' '19.02'
Set-TextValue $ws 'E28' '  +2.31%  '
Set-TextValue $ws 'D29' '2.098'
Set-TextValue $ws 'E29' '  +4.09%  '
Set-TextValue $ws 'D30' '1.374'
Set-TextValue $ws 'E30' '  -1.92%  '
Set-TextValue $ws 'D31' '1.519'
Set-TextValue $ws 'E31' '  +0.24%  '
Set-TextValue $ws 'D32' '4.321'
Set-TextValue $ws 'E32' '  +1.69%  '
Set-TextValue $ws 'D33' '0.05686'
Set-TextValue $ws 'E33' '  +9.12%  '
Set-TextValue $ws 'D34' '4.100'
Set-TextValue $ws 'E34' '  +1.04%  '
Set-TextValue $ws 'D35' '1.297'
Set-TextValue $ws 'E35' '  +4.22%  '
Set-TextValue $ws 'D36' '0.7374'
Set-TextValue $ws 'E36' '  +1.97%  '
Set-TextValue $ws 'D37' '2.723'
Set-TextValue $ws 'E37' '  +0.31%  '
Set-TextValue $ws 'D38' '0.01923'
Set-TextValue $ws 'E38' '  +0.55%  '
Set-TextValue $ws 'D39' '2.796'
Set-TextValue $ws 'E39' '  +0.98%  '
Set-TextValue $ws 'D40' '0.4459'
Set-TextValue $ws 'E40' '  +1.55%  '
Set-TextValue $ws 'D41' '72.68'
Set-TextValue $ws 'E41' '  +1.73%  '
Set-TextValue $ws 'D42' '6.025'
Set-TextValue $ws 'E42' '  -1.66%  '
Set-TextValue $ws 'D43' '0.8444'
Set-TextValue $ws 'E43' '  +1.94%  '
Set-TextValue $ws 'D44' '1.908'
Set-TextValue $ws 'E44' '  +1.67%  '
Set-TextValue $ws 'D45' '1.001'
Set-TextValue $ws 'E45' '  -0.05%  '
Set-TextValue $ws 'D46' '7.648'
Set-TextValue $ws 'E46' '  +0.97%  '
Set-TextValue $ws 'D47' '101.32'
Set-TextValue $ws 'E47' '  +1.59%  '
Set-TextValue $ws 'D48' '9.817'
Set-TextValue $ws 'E48' '  +0.87%  '
Set-TextValue $ws 'D49' '991.77'
Set-TextValue $ws 'E49' '  +9.83%  '
Set-TextValue $ws 'D50' '2.061.95'
Set-TextValue $ws 'E50' '  +0.77%  '
Set-TextValue $ws 'D51' '36.49'
Set-TextValue $ws 'E51' '  +1.65%  '
